$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (peilvak, wl_node_up,
# wl_node_down) from A:C to B:D.
$ws.Range("A1").EntireColumn.Insert()

# Add the new "id" header and row id values 1, 2, 3.
$ws.Range("A1").Value = "id"
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3

# Match the new active selection recorded in the workbook (cell A5).
$ws.Range("A5").Select()
